$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new timesheet entry as row 35, following the same pattern as
# the existing rows: a date in column A, hours worked in column B, and a
# running total formula in column C.

# Copy the date cell's number format from the previous row (A34) so the
# new date (A35) renders the same way, then set its value.
$ws.Range("A34").Copy() | Out-Null
$ws.Range("A35").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A35").Value = 45345

# Hours worked for the new entry.
$ws.Range("B35").Value = 0.5

# Running total: previous total (C34) plus the new hours (B35).
$ws.Range("C35").Formula = "=C34+B35"

$excel.CutCopyMode = 0

# Reflect the newly-entered cell as the active selection, matching how a
# user would have just finished typing the new row.
$ws.Range("B35").Select() | Out-Null

$wb.Save()
